$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 388
$wsExhibit.Range("F4").Value = 0
$wsExhibit.Range("F5").Value = 0
$wsExhibit.Range("F7").Value = 0
$wsExhibit.Range("F8").Value = 0
$wsExhibit.Range("F10").Value = 486

# Sheet "全部类型" (sheet4) updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 0
$wsAll.Range("F4").Value = 1626
$wsAll.Range("F5").Value = 0
$wsAll.Range("F6").Value = 0
$wsAll.Range("F7").Value = 416
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 486
